# Apply updated completeness_rank (column E) values as computed in the
# "calculate missing and non-zero stats v2" recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "E16"  = 142
    "E17"  = 143
    "E31"  = 127
    "E32"  = 128
    "E53"  = 106
    "E54"  = 105
    "E59"  = 100
    "E60"  = 99
    "E63"  = 95
    "E64"  = 96
    "E69"  = 90
    "E70"  = 89
    "E76"  = 82
    "E77"  = 83
    "E78"  = 80
    "E80"  = 79
    "E81"  = 78
    "E82"  = 75
    "E83"  = 76
    "E84"  = 77
    "E86"  = 73
    "E87"  = 71
    "E88"  = 74
    "E89"  = 67
    "E90"  = 69
    "E91"  = 68
    "E95"  = 63
    "E96"  = 64
    "E97"  = 62
    "E98"  = 59
    "E99"  = 60
    "E100" = 61
    "E101" = 57
    "E102" = 58
    "E103" = 56
    "E105" = 54
    "E110" = 45
    "E111" = 47
    "E112" = 46
    "E113" = 49
    "E114" = 48
    "E115" = 43
    "E116" = 44
    "E117" = 38
    "E118" = 42
    "E119" = 41
    "E120" = 39
    "E121" = 40
    "E122" = 37
    "E123" = 35
    "E124" = 34
    "E125" = 36
    "E126" = 33
    "E127" = 32
    "E128" = 31
    "E129" = 29
    "E130" = 27
    "E131" = 28
    "E132" = 30
    "E133" = 26
    "E134" = 24
    "E136" = 25
    "E140" = 17
    "E141" = 18
    "E142" = 16
    "E143" = 19
    "E148" = 11
    "E149" = 10
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
